$d = $word.ActiveDocument

# --------------------------------------------------------------------
# 1) Paragraph "Nome: Priscila Cristina RA.:820147927":
#    the three runs (split by a pair of gramStart/gramEnd proofErr
#    markers) collapse into a single run with the same visible text.
#    Replacing the full phrase with itself via Find/Replace makes Word
#    re-emit it as one run and drops the now-unneeded proofErr marks.
# --------------------------------------------------------------------
$d.Content.Find.Execute("Nome: Priscila Cristina RA.:820147927", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Nome: Priscila Cristina RA.:820147927", 2) | Out-Null

# --------------------------------------------------------------------
# 2) Paragraph "Nome: Vinicius Kulik Gavioli ...": the leading
#    "Nome: " and "Vinicius " runs merge into a single run. Again a
#    self-replace of that leading phrase is enough to coalesce them.
# --------------------------------------------------------------------
$d.Content.Find.Execute("Nome: Vinicius ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Nome: Vinicius ", 2) | Out-Null

# --------------------------------------------------------------------
# 3) Same paragraph: append a new "]" run right after the
#    " R.A.:819151742" text (still inside this paragraph, and - at
#    this point - before the "_GoBack" bookmark, which we are about to
#    relocate in step 4).
# --------------------------------------------------------------------
$pVinicius = $d.Paragraphs(4)
$rVinicius = $pVinicius.Range
$viniciusTextLen = $rVinicius.End - $rVinicius.Start - 1
$viniciusEndPos = $rVinicius.Start + $viniciusTextLen
$viniciusEndRange = $d.Range($viniciusEndPos, $viniciusEndPos)
$viniciusEndRange.InsertAfter("]")

# --------------------------------------------------------------------
# 4) Insert a brand-new paragraph after the Vinicius paragraph holding
#    "Nome: Marcelo Vinicius Martins da Silva R.A: 820134048", and move
#    the "_GoBack" bookmark (previously at the end of the Vinicius
#    paragraph) onto the end of this new paragraph.
# --------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$pVinicius = $d.Paragraphs(4)
$pVinicius.Range.InsertParagraphAfter()
$pMarcelo = $d.Paragraphs(5)

$marceloText = "Nome: Marcelo Vinicius Martins da Silva R.A: 820134048"
# Type the real text plus one throw-away trailing character so the
# bookmark's insertion point below is an interior position of the
# paragraph rather than sitting right on the paragraph mark.
$pMarcelo.Range.Text = $marceloText + "X"
$rMarcelo = $pMarcelo.Range
$bookmarkPos = $rMarcelo.Start + $marceloText.Length
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# Remove the throw-away character now that the bookmark is anchored
# right after the real text.
$dummyRange = $d.Range($bookmarkPos, $bookmarkPos + 1)
$dummyRange.Delete()
